# "updated the class as done"
#
# Tasks sheet ("Tasks", sheet1):
#   - B10 and B11 get marked "Done" (same shared string/style already used
#     by B18/B19 for the "EAICommon"/"EJBServer" rows).
#   - The user's cursor ends up on B12, having scrolled the window down so
#     row 13 is pinned to the top.
#
# "Serviços In Out" sheet (sheet2):
#   - Row 19 (an otherwise empty spacer row) had an explicit 14.45pt custom
#     height; it goes back to the sheet's default row height.

$wb = $excel.ActiveWorkbook

# --- "Serviços In Out": drop the custom height on row 19 ---------------
$wsServices = $wb.Worksheets.Item(2)
[void]$wsServices.Activate()
[void]$wsServices.Rows(19).AutoFit()

# --- "Tasks": fill in the new "Done" cells and update the view ---------
$wsTasks = $wb.Worksheets.Item(1)
[void]$wsTasks.Activate()

$wsTasks.Range("B10").Value = "Done"
$wsTasks.Range("B11").Value = "Done"

# Scroll so row 13 is at the top of the window, then land the selection on
# B12, matching where the editor left the cursor.
[void]$wsTasks.Range("A13").Select()
$excel.ActiveWindow.ScrollRow = 13
$excel.ActiveWindow.ScrollColumn = 1
[void]$wsTasks.Range("B12").Select()
